$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "300.79"
    "E2" = "-0.73%"
    "D3" = "31.37"
    "E3" = "-2.19%"
    "D4" = "5.133"
    "E4" = "-2.57%"
    "D5" = "0.07369"
    "E5" = "-1.45%"
    "D6" = "2.299"
    "E6" = "52.53%"
    "D7" = "7.963"
    "E7" = "1.46%"
    "D8" = "3.793"
    "E8" = "-0.48%"
    "D9" = "0.9196"
    "E9" = "-0.16%"
    "D10" = "0.1716"
    "E10" = "2.01%"
    "D11" = "0.07652"
    "E11" = "-4.38%"
    "D12" = "0.08125"
    "E12" = "1.52%"
    "D13" = "0.03016"
    "E13" = "-0.41%"
    "E14" = "0.43%"
    "D15" = "0.001495"
    "E15" = "0.04%"
    "D16" = "0.006194"
    "E16" = "-2.39%"
    "E17" = "0.04%"
    "D18" = "2.227"
    "E18" = "-0.23%"
    "E19" = "-0.54%"
    "E20" = "-0.38%"
    "D21" = "4.654"
    "E21" = "3.43%"
    "D22" = "0.04633"
    "E22" = "0.63%"
    "E23" = "-3.20%"
    "D24" = "0.001225"
    "E24" = "0.69%"
    "D25" = "0.004486"
    "E26" = "-7.25%"
    "E27" = "5.52%"
    "D39" = "0.01735"
    "E39" = "1.63%"
    "D40" = "0.04522"
    "E40" = "0.97%"
    "D41" = "0.007221"
    "E41" = "1.06%"
    "D42" = "0.1346"
    "E42" = "-0.41%"
    "D43" = "0.002209"
    "E43" = "-1.04%"
    "D44" = "0.01071"
    "E44" = "-16.31%"
    "D45" = "0.00006272"
    "E45" = "1.79%"
    "D47" = "1.928"
    "E47" = "3.36%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}
